# Update "想去人数" (interested count) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 968
$wsExhibit.Range("F10").Value = 489
$wsExhibit.Range("F14").Value = 4547
$wsExhibit.Range("F17").Value = 2834
$wsExhibit.Range("F20").Value = 1143
$wsExhibit.Range("F21").Value = 3844
$wsExhibit.Range("F24").Value = 1544
$wsExhibit.Range("F26").Value = 2513
$wsExhibit.Range("F34").Value = 270
$wsExhibit.Range("F35").Value = 47
$wsExhibit.Range("F36").Value = 94
$wsExhibit.Range("F37").Value = 1472
$wsExhibit.Range("F41").Value = 15
$wsExhibit.Range("F43").Value = 138
$wsExhibit.Range("F45").Value = 327
$wsExhibit.Range("F49").Value = 93

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 968
$wsAll.Range("F8").Value = 489
$wsAll.Range("F10").Value = 4548
$wsAll.Range("F14").Value = 2834
$wsAll.Range("F16").Value = 1143
$wsAll.Range("F17").Value = 3844
$wsAll.Range("F20").Value = 1544
$wsAll.Range("F23").Value = 2513
$wsAll.Range("F34").Value = 270
$wsAll.Range("F35").Value = 1472
$wsAll.Range("F40").Value = 15
$wsAll.Range("F43").Value = 138
$wsAll.Range("F45").Value = 327
$wsAll.Range("F49").Value = 93
